# Add data for 2021-09-30
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2021-09-22"

# Row 3 (January) - 2021 column (arrest_made/no_arrest_made/arrest_rate) updated
$ws.Range("T3").Value = 15
$ws.Range("U3").Value = 201
$ws.Range("V3").Value = 0.0694

# Row 11 (September) - label and 2017/2018/2019/2020/2021 values updated
$ws.Range("A11").Value = "September (through 09-22)"
$ws.Range("I11").Value = 47
$ws.Range("J11").Value = 0.0784
$ws.Range("K11").Value = 4
$ws.Range("L11").Value = 38
$ws.Range("M11").Value = 0.09520000000000001
$ws.Range("O11").Value = 51
$ws.Range("P11").Value = 0.0727
$ws.Range("R11").Value = 81
$ws.Range("S11").Value = 0.0357
$ws.Range("U11").Value = 133

# Row 12 (Total)
$ws.Range("I12").Value = 553
$ws.Range("J12").Value = 0.0814
$ws.Range("K12").Value = 61
$ws.Range("L12").Value = 471
$ws.Range("M12").Value = 0.1147
$ws.Range("O12").Value = 364
$ws.Range("P12").Value = 0.099
$ws.Range("R12").Value = 817
$ws.Range("S12").Value = 0.0588
$ws.Range("T12").Value = 73
$ws.Range("U12").Value = 1130
$ws.Range("V12").Value = 0.0607
